$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 1.742479666666667
$arr[0,1] = 5.227439
$arr[0,2] = 0.1294580684571358
$arr[0,3] = 0.1294580684571358
$ws.Range("G2:J2").Value = $arr
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 2.989238333333333
$arr[0,1] = 8.967715
$arr[0,2] = 0.2714069173313896
$arr[0,3] = 0.2714069173313896
$arr[0,4] = 5.208687014653889
$arr[0,5] = 46.878183131885
$arr[0,6] = 0.03513581528362723
$arr[0,7] = 0.03513581528362723
$ws.Range("M2:T2").Value = $arr

# Row 3
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 1.742479666666667
$arr[0,1] = 5.227439
$arr[0,2] = 0.1294580684571358
$arr[0,3] = 0.1294580684571358
$ws.Range("G3:J3").Value = $arr
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.4642544095794917
$arr[0,1] = 0.4642544095794917
$arr[0,2] = 8.909706275908666
$arr[0,3] = 80.18735648317801
$arr[0,4] = 0.06010147913686901
$arr[0,5] = 0.06010147913686899
$ws.Range("O3:T3").Value = $arr

# Row 4
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 1.742479666666667
$arr[0,1] = 5.227439
$arr[0,2] = 0.1294580684571358
$arr[0,3] = 0.1294580684571358
$ws.Range("G4:J4").Value = $arr
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.2643386730891187
$arr[0,1] = 0.2643386730891186
$arr[0,2] = 5.073037295910112
$arr[0,3] = 45.657335663191
$arr[0,4] = 0.03422077403663957
$arr[0,5] = 0.03422077403663955
$ws.Range("O4:T4").Value = $arr

# Row 5
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.5665027357143181
$arr[0,1] = 0.5665027357143181
$ws.Range("I5:J5").Value = $arr
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 2.989238333333333
$arr[0,1] = 8.967715
$arr[0,2] = 0.2714069173313896
$arr[0,3] = 0.2714069173313896
$arr[0,4] = 22.79298214817777
$arr[0,5] = 205.1368393336
$arr[0,6] = 0.153752761160022
$arr[0,7] = 0.153752761160022
$ws.Range("M5:T5").Value = $arr

# Row 6
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.5665027357143181
$arr[0,1] = 0.5665027357143181
$ws.Range("I6:J6").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.4642544095794917
$arr[0,1] = 0.4642544095794917
$ws.Range("O6:P6").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.2630013930942176
$arr[0,1] = 0.2630013930942176
$ws.Range("S6:T6").Value = $arr

# Row 7
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.5665027357143181
$arr[0,1] = 0.5665027357143181
$ws.Range("I7:J7").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.2643386730891187
$arr[0,1] = 0.2643386730891186
$ws.Range("O7:P7").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.1497485814600785
$arr[0,1] = 0.1497485814600785
$ws.Range("S7:T7").Value = $arr

# Row 8
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.3040391958285462
$arr[0,1] = 0.3040391958285461
$ws.Range("I8:J8").Value = $arr
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 2.989238333333333
$arr[0,1] = 8.967715
$arr[0,2] = 0.2714069173313896
$arr[0,3] = 0.2714069173313896
$arr[0,4] = 12.23287995975556
$arr[0,5] = 110.0959196378
$arr[0,6] = 0.0825183408877404
$arr[0,7] = 0.08251834088774039
$ws.Range("M8:T8").Value = $arr

# Row 9
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.3040391958285462
$arr[0,1] = 0.3040391958285461
$ws.Range("I9:J9").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.4642544095794917
$arr[0,1] = 0.4642544095794917
$ws.Range("O9:P9").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.1411515373484052
$arr[0,1] = 0.1411515373484051
$ws.Range("S9:T9").Value = $arr

# Row 10
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.3040391958285462
$arr[0,1] = 0.3040391958285461
$ws.Range("I10:J10").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.2643386730891187
$arr[0,1] = 0.2643386730891186
$ws.Range("O10:P10").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.0803693175924006
$arr[0,1] = 0.08036931759240057
$ws.Range("S10:T10").Value = $arr

